# Applies the "Added a few more slots" edit to the Double Dinosaur Deluxe
# review document:
#   1. Insert a new "Meta description: ..." paragraph right after the
#      Heading1 title at the top of the document (bold label, plain value).
#   2. Remove the duplicate bold "Play Double Dinosaur Deluxe for Free -
#      Slot Review" paragraph that was sitting near the end of the doc.
#   3. Replace the text of the trailing italic paragraph (formerly the meta
#      description) with the new AI image-generation prompt, keeping its
#      italic formatting intact.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: swap the trailing italic paragraph's text for the new image
# prompt *before* we duplicate its old text elsewhere, so the Find call
# below can't accidentally match the wrong (freshly inserted) paragraph.
# ---------------------------------------------------------------------
$oldMetaText = "Read our review of Double Dinosaur Deluxe online slot game and play for free. Features, gameplay, graphics, winning potential, and more."
$newImagePrompt = "Create a feature image for Double Dinosaur Deluxe that showcases the game's fun and adventurous theme. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The warrior should be surrounded by dinosaurs and be holding a map or a binocular to show their adventurous spirit. The background should be of a prehistoric landscape with a smoking volcano in the distance to add to the game's suspenseful atmosphere. Make sure to use bright colors and bold lines to make the image pop and grab the attention of potential players."

$italicTarget = $d.Content
$foundItalic = $italicTarget.Find.Execute($oldMetaText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundItalic) {
    # Setting .Text directly (rather than Find's Replace argument) keeps
    # straight apostrophes instead of letting autocorrect curl them, and
    # preserves the run's existing <w:rPr> (the italic formatting).
    $italicTarget.Text = $newImagePrompt
}

# ---------------------------------------------------------------------
# Step 2: delete the old, now-redundant bold "Play Double Dinosaur
# Deluxe for Free - Slot Review" paragraph near the end of the document
# (the Heading1 at the very top of the doc keeps the same text, so match
# on the Normal-styled copy specifically).
# ---------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.ParagraphStyle.NameLocal -eq "Normal" -and `
        $p.Range.Text -match "Play Double Dinosaur Deluxe for Free - Slot Review") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# Step 3: insert the new "Meta description: ..." paragraph right after
# the Heading1 title paragraph at the top.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"

$metaRange = $metaPara.Range
$metaStart = $metaRange.Start
$label = "Meta description"
$rest = ": Read our review of Double Dinosaur Deluxe online slot game and play for free. Features, gameplay, graphics, winning potential, and more."
$metaRange.Text = $label + $rest

# Bold just the "Meta description" label, leaving the rest of the
# sentence in regular formatting.
$labelRange = $d.Range($metaStart, $metaStart + $label.Length)
$labelRange.Bold = 1

Write-Output "ok"
